$d = $word.ActiveDocument

# The paragraph currently contains a single run holding the inline drawing,
# followed by the (hidden) "_GoBack" bookmark. We need to prepend the text
# "Class Diagram: by Jacob Field" (split into two runs with the _GoBack
# bookmark sitting between "Class" and " Diagram: by Jacob Field"), mark
# every run in the paragraph as <w:noProof/>, and relocate the bookmark from
# after the drawing to between the two new text runs.

$firstPara = $d.Paragraphs(1)

# Insert a brand-new empty paragraph ahead of the drawing paragraph and fill
# it with the caption text. Using InsertParagraphBefore + Range.Text (rather
# than Range.InsertBefore/InsertAfter directly against the drawing's run)
# keeps the drawing run intact instead of having its contents clobbered.
$firstPara.Range.InsertParagraphBefore()
$textPara = $d.Paragraphs(1)
$textPara.Range.Text = "Class Diagram: by Jacob Field"

# Merge the new text paragraph into the drawing paragraph by deleting the
# paragraph mark between them.
$mergeEnd = $d.Paragraphs(1).Range.End
$d.Range($mergeEnd - 1, $mergeEnd).Delete()

# Split "Class Diagram: by Jacob Field" into "Class" and
# " Diagram: by Jacob Field" runs. Toggling a character formatting property
# on/off over just the "Class" span forces the engine to break the run
# without leaving any residual formatting behind.
$splitPoint = 5
$classRange = $d.Range(0, $splitPoint)
$classRange.Bold = 1
$classRange.Bold = 0

# Mark every run up to (and including) the drawing as noProof.
$d.Range(0, $splitPoint).NoProofing = $true
$d.Range($splitPoint, 30).NoProofing = $true

# Relocate the "_GoBack" bookmark from right after the drawing to right
# between the two text runs ("Class" | " Diagram: by Jacob Field").
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()
$d.Bookmarks.Add("_GoBack", $d.Range($splitPoint, $splitPoint)) | Out-Null
